$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row additions ---
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Qualification"

# --- Column A: UserName values ---
$names = @("dbhatt","raj","mohita","Balram","Nellam","Sonu","Hema","Namita")
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
}

# --- Column B: Password values ---
$ws.Range("B2").Value = "Tesy@123"
for ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = "Test@123"
}

# --- Column C: Email values ---
$ws.Range("C2").Value = "damomca@gmail.com"
for ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = "user@user.com"
}

# --- Column D: Qualification values ---
$quals = @("MCA","BCA","BTECH","BE","MBA","CS","BCOM","BTECH")
for ($i = 0; $i -lt $quals.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 4).Value = $quals[$i]
}

# --- Hyperlinks (auto-generated by Excel for "x@y"-shaped text) ---
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Tesy@123")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("B4:B9"), "mailto:Test@123", "", "", "Test@123")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:damomca@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:user@user.com")
$ws.Hyperlinks.Add($ws.Range("C4:C9"), "mailto:user@user.com", "", "", "user@user.com")

# Apply the Hyperlink style to every linked cell (Add() only styles the first cell of a range)
$ws.Range("B2:B9").Style = "Hyperlink"
$ws.Range("C2:C9").Style = "Hyperlink"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 13.1
$ws.Columns.Item(2).ColumnWidth = 14.1
$ws.Columns.Item(3).ColumnWidth = 29.5
$ws.Columns.Item(4).ColumnWidth = 12.5

# --- Selection ---
[void]$ws.Range("E7").Select()
